# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded on the localization
# status sheets (Overview, zh-cn, de-de) to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-14 17:31:33"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-14 17:31:25"
$wsZhCn.Range("K2").Value = "2016-08-14 17:32:07"

# de-de sheet: Correspond Handoff / Handback datetimes for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-14 17:31:33"
$wsDeDe.Range("K2").Value = "2016-08-14 17:32:17"
